$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting refreshed crypto price/volume/hour data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '308.77'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.19%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '12'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.44%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '12'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.114'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.69%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '12'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07631'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.55%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '12'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.252'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.33%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '12'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.608'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.72%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '12'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.517'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '4.01%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '12'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9044'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.42%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '12'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1115'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '9.74%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '12'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1794'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '3.87%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '12'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09107'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.15%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '12'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04156'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-6.62%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '12'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1051'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.21%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '12'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001254'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.08%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '12'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005790'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.31%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '12'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.338'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.49%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '12'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.23%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '12'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.835'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-4.02%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '12'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1364'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.90%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '12'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2703'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.39%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '12'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04053'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.09%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '12'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001244'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.23%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '12'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004097'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.77%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '12'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001300'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.14%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '12'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003742'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-95.01%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '12'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '12'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '12'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '12'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '12'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '12'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '12'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '12'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '12'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '12'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '12'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '12'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02394'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '2.13%'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '12'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05231'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.46%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '12'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007772'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.91%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '12'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1300'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.41%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '12'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007042'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '13.57%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '12'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001950'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.49%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '12'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007725'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-6.00%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '12'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3345'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.03%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '12'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006918'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '6.21%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '12'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.16%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '12'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.04919'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '1,372.24%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '12'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '12'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.16%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '12'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.16%'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '12'
